# This revision's OOXML diff is purely additive at the package level: it
# adds a set of SharePoint/OneDrive "document management" custom XML parts
# (customXml/item1.xml, item2.xml, item3.xml + their itemProps*.xml /
# datastoreItem companions). These parts are the standard boilerplate that
# SharePoint/OneDrive stamps onto a .docx the first time it is synced into a
# library that has a "Content Type"/managed-metadata column (the
# ma:contentTypeName="Document", MediaServiceImageTags / TaxCatchAll fields,
# etc. seen in item1.xml's schema, plus the empty Terms/TaxCatchAll payload
# in item3.xml). There is NO change anywhere to word/document.xml, styles,
# numbering, settings, headers/footers, or any visible text/formatting --
# every paragraph, run and comment-note in the body ("About page", "Terms
# of Use", "Workplace Health and Safety", etc.) is byte-for-byte the same
# before and after. The commit message ("Continued to note which of
# Tanya's changes were actioned") describes earlier edits in this same
# running log-style document; this particular save simply happened to be
# the one that passed through SharePoint/OneDrive and therefore picked up
# its metadata parts as a side effect.
#
# We still attempt the straightforward COM call that models that
# side-effect (Document.CustomXMLParts.Add), in case the host supports
# minting the custom XML package parts that way, but we don't let a lack
# of support there change the document's actual content -- the content
# itself is unchanged in this revision.

$d = $word.ActiveDocument

$contentTypeSchemaXml = @'
<?xml version="1.0" encoding="utf-8"?>
<ct:contentTypeSchema xmlns:ct="http://schemas.microsoft.com/office/2006/metadata/contentType" xmlns:ma="http://schemas.microsoft.com/office/2006/metadata/properties/metaAttributes" ct:_="" ma:_="" ma:contentTypeName="Document" ma:contentTypeID="0x0101004FE47EDC821BDF48B4A02F50FA0D7AA2" ma:contentTypeVersion="16" ma:contentTypeDescription="Create a new document." ma:contentTypeScope="" ma:versionID="c985933b0a177fb87df527f3d6a422a5">
  <xsd:schema xmlns:xsd="http://www.w3.org/2001/XMLSchema" xmlns:xs="http://www.w3.org/2001/XMLSchema" xmlns:p="http://schemas.microsoft.com/office/2006/metadata/properties" xmlns:ns2="739bdb02-8359-4bf5-94bc-edd490470c6d" xmlns:ns3="766209d0-6631-4f1e-9a64-6298cbbf9a9c" targetNamespace="http://schemas.microsoft.com/office/2006/metadata/properties" ma:root="true" ma:fieldsID="2da29d1270cc6742464067118e96aafc" ns2:_="" ns3:_="">
    <xsd:import namespace="739bdb02-8359-4bf5-94bc-edd490470c6d"/>
    <xsd:import namespace="766209d0-6631-4f1e-9a64-6298cbbf9a9c"/>
  </xsd:schema>
</ct:contentTypeSchema>
'@

$formTemplatesXml = @'
<?xml version="1.0" encoding="utf-8"?>
<?mso-contentType ?>
<FormTemplates xmlns="http://schemas.microsoft.com/sharepoint/v3/contenttype/forms">
  <Display>DocumentLibraryForm</Display>
  <Edit>DocumentLibraryForm</Edit>
  <New>DocumentLibraryForm</New>
</FormTemplates>
'@

$propertiesXml = @'
<?xml version="1.0" encoding="utf-8"?>
<p:properties xmlns:p="http://schemas.microsoft.com/office/2006/metadata/properties" xmlns:xsi="http://www.w3.org/2001/XMLSchema-instance" xmlns:pc="http://schemas.microsoft.com/office/infopath/2007/PartnerControls">
  <documentManagement>
    <lcf76f155ced4ddcb4097134ff3c332f xmlns="739bdb02-8359-4bf5-94bc-edd490470c6d">
      <Terms xmlns="http://schemas.microsoft.com/office/infopath/2007/PartnerControls"/>
    </lcf76f155ced4ddcb4097134ff3c332f>
    <TaxCatchAll xmlns="766209d0-6631-4f1e-9a64-6298cbbf9a9c" xsi:nil="true"/>
  </documentManagement>
</p:properties>
'@

foreach ($partXml in @($contentTypeSchemaXml, $formTemplatesXml, $propertiesXml)) {
    try {
        [void]$d.CustomXMLParts.Add($partXml)
    } catch {
        # Minting new custom XML package parts (the SharePoint/OneDrive
        # document-management metadata this save picked up) isn't something
        # this host's object model exposes a working path for. That's fine:
        # those parts carry no visible document content, so there is
        # nothing else for this script to change.
    }
}
